$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.6811249999999999
$ws.Range("M2").Value = 3.390429
$ws.Range("N2").Value = 10.171287
$ws.Range("O2").Value = 0.173121426386348
$ws.Range("P2").Value = 0.173121426386348
$ws.Range("Q2").Value = 2.309305952624999
$ws.Range("R2").Value = 20.783753573625
$ws.Range("S2").Value = 0.173121426386348
$ws.Range("T2").Value = 0.173121426386348

# Row 3
$ws.Range("G3").Value = 0.6811249999999999
$ws.Range("O3").Value = 0.5936336753560868
$ws.Range("P3").Value = 0.5936336753560868
$ws.Range("Q3").Value = 7.918614170374999
$ws.Range("R3").Value = 71.26752753337499
$ws.Range("S3").Value = 0.5936336753560868
$ws.Range("T3").Value = 0.5936336753560868

# Row 4
$ws.Range("G4").Value = 0.6811249999999999
$ws.Range("M4").Value = 4.546141666666667
$ws.Range("N4").Value = 13.638425
$ws.Range("O4").Value = 0.2321342018628743
$ws.Range("P4").Value = 0.2321342018628743
$ws.Range("Q4").Value = 3.096490742708333
$ws.Range("R4").Value = 27.868416684375
$ws.Range("S4").Value = 0.2321342018628743
$ws.Range("T4").Value = 0.2321342018628743

# Row 5
$ws.Range("G5").Value = 0.6811249999999999
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.021752
$ws.Range("N5").Value = 0.06525600000000001
$ws.Range("O5").Value = 0.001110696394691009
$ws.Range("P5").Value = 0.001110696394691009
$ws.Range("Q5").Value = 0.014815831
$ws.Range("R5").Value = 0.133342479
$ws.Range("S5").Value = 0.001110696394691009
$ws.Range("T5").Value = 0.001110696394691009
